# Applies the "pushing latest status. Transform app and its output added" commit
# to the "To do List" workbook:
#   - Adds a Status of "Complete" to rows 3, 11 and 28 (column G)
#   - Adds a Status of "Complete - Pending review" to row 20 (column G)
#   - Adds a Note to row 28 (column E): "Bad values in data were preventing
#     auto scale even after Forced Zero was toggled"
#   - Shrinks row 22's height (it no longer needs as much space)
#   - Scrolls/freezes the sheet further down and updates the active selection
#     to reflect the latest place being worked on (around row 24)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status column (G) updates -------------------------------------------
$ws.Range("G3").Value  = "Complete"
$ws.Range("G11").Value = "Complete"
$ws.Range("G20").Value = "Complete - Pending review"
$ws.Range("G28").Value = "Complete"

# --- Notes column (E) update for row 28 -----------------------------------
$ws.Range("E28").Value = "Bad values in data were preventing auto scale even after Forced Zero was toggled"

# Make sure the newly populated Status cells keep the same look (wrapped
# text) as the rest of column G.
$ws.Range("G3").Style  = $ws.Range("G2").Style
$ws.Range("G11").Style = $ws.Range("G2").Style
$ws.Range("G20").Style = $ws.Range("G2").Style
$ws.Range("G28").Style = $ws.Range("G2").Style
$ws.Range("E28").Style = $ws.Range("E2").Style

# --- Row height: item 4.1 no longer needs the extra room ------------------
$ws.Rows.Item(22).RowHeight = 45

# --- Update frozen pane / scroll position & active selection --------------
$ws.Range("G24").Select()
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A25").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("G24").Select()
